# This script collapses several runs of "word <space> word <space> ..." text
# in the cover-page paragraphs (and the TOC heading) into a single run each,
# matching the upstream regenerated-docx output. The wording/content itself
# is unchanged except the report date, which moves from
# "November 18, 2024" to "November 27, 2024".
#
# Word's Find/Replace merges the matched span into one run (with the
# formatting of the first run in the span), which is exactly the
# transformation the diff shows, so we simply replace each already-correct
# phrase with itself (or, for the date, with the new date).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute(
        $find,      # FindText
        $true,      # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $replace,   # ReplaceWith
        2           # Replace (wdReplaceOne)
    ) | Out-Null
}

# Title paragraph
Replace-Text "Analytics for the Australian Grains Industry - Curtin University (AAGI-CU)" "Analytics for the Australian Grains Industry - Curtin University (AAGI-CU)"
Replace-Text "Technical Report Series: 123" "Technical Report Series: 123"

# Subtitle paragraph
Replace-Text "Descriptive title for report" "Descriptive title for report"
Replace-Text "Report for AAA–BBB" "Report for AAA–BBB"

# Author paragraphs
Replace-Text "Prepared for: collab_partner (collab_partner@email.com.au)" "Prepared for: collab_partner (collab_partner@email.com.au)"
Replace-Text "Prepared by: Your.Name" "Prepared by: Your.Name"
Replace-Text "Project Lead: Curtin University – Prof Mark Gibberd, Dr Julia Easton, Prof Adam Sparks" "Project Lead: Curtin University – Prof Mark Gibberd, Dr Julia Easton, Prof Adam Sparks"
Replace-Text "email: cbada@curtin.edu.au" "email: cbada@curtin.edu.au"

# Date paragraph - content actually changes here
Replace-Text "November 18, 2024" "November 27, 2024"

# Table of contents heading
# This paragraph lives inside a w:sdt (TOC building-block content control).
# Find/Execute does not traverse into sdt content in this runtime, so locate
# the paragraph by style name and collapse its runs by assigning Range.Text
# directly instead. A no-op assignment (same text in, same text out) doesn't
# trigger the run-merge, so nudge it through a distinct intermediate value
# first, then set the real text.
$tocHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Style.NameLocal -eq "TOC Heading") {
        $tocHeadingIndex = $i
        break
    }
}

if ($tocHeadingIndex -gt 0) {
    $tocPara = $d.Paragraphs.Item($tocHeadingIndex)
    $tocRange = $d.Range($tocPara.Range.Start, $tocPara.Range.End - 1)
    $tocRange.Text = "Table of contents#"

    $tocPara2 = $d.Paragraphs.Item($tocHeadingIndex)
    $tocRange2 = $d.Range($tocPara2.Range.Start, $tocPara2.Range.End - 1)
    $tocRange2.Text = "Table of contents"
}
